$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.491.66"
$ws.Range("E2").Value = "  +3.60%  "
$ws.Range("D3").Value = "1.588.34"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("E4").Value = "  +0.99%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("E7").Value = "  +1.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.29"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.52%  "
$ws.Range("E9").Value = "  +0.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0600"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0886"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.73%  "
$ws.Range("D12").Value = "1.817.37"
$ws.Range("E12").Value = "  +1.22%  "
$ws.Range("D13").Value = "1.592.77"
$ws.Range("E13").Value = "  +1.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.528"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.76%  "
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").Value = "28.508.74"
$ws.Range("E16").Value = "  +3.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.11"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.99%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.94%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0707"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.80%  "
$ws.Range("E21").Value = "  +0.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.08%  "
$ws.Range("E24").Value = "  +1.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.63%  "
$ws.Range("E27").Value = "  -0.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.106"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.66%  "
$ws.Range("E29").Value = "  +0.99%  "
$ws.Range("E30").Value = "  -0.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0471"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.95%  "
$ws.Range("D34").Value = "1.403.96"
$ws.Range("E34").Value = "  -3.64%  "
$ws.Range("E35").Value = "  -1.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -10.24%  "
$ws.Range("E37").Value = "  +1.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.61"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +10.26%  "
$ws.Range("E39").Value = "  -0.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.542"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.812"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("E42").Value = "  +0.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.63"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.43%  "
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.984"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.02"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.95%  "
$ws.Range("D47").Value = "1.727.27"
$ws.Range("E47").Value = "  +1.19%  "
$ws.Range("B48").Value = "mCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.39%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "87.28"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.46%  "
$ws.Range("D50").Value = "0.0₆0104"
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0522"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.59%  "
